$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated coin price/volume data as captured in the latest GitHub Actions run.
# Columns B (Coin) and C (Link) are plain text; columns D (Price) and E (Volume(1h))
# hold numeric-looking strings, so their NumberFormat is forced to Text ("@") before
# assignment to prevent Excel from auto-converting them into numeric/percentage values.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '312.96'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '1.44%'

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '39.47'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '2.62%'

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.152'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '1.11%'

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.08182'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '0.71%'

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.976'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-0.14%'

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '8.139'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '2.54%'

# Row 8
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9275'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '0.01%'

# Row 9
$ws.Range("B9").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C9").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1387'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-3.70%'

# Row 10
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1933'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-1.26%'

# Row 11
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.09048'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-0.34%'

# Row 12
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03504'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-0.09%'

# Row 13
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09792'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-0.19%'

# Row 14
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001391'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-0.80%'

# Row 15
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.006131'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '0.22%'

# Row 16
$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.676'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '0.49%'

# Row 17
$ws.Range("B17").Value = 'GateToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.229'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '0.53%'

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3460'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '0.83%'

# Row 20
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '1.23%'

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.637'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-3.43%'

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04367'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-1.14%'

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001236'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '1.61%'

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004812'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-0.47%'

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0001298'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-0.33%'

# Row 27
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '-10.38%'

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02168'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '3.06%'

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.05197'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '0.91%'

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007456'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-0.32%'

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.009723'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-4.03%'

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1376'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '1.37%'

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.002115'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-1.39%'

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.009858'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '6.98%'

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006357'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '1.22%'

# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-0.57%'

# Row 48
$ws.Range("B48").Value = 'CoinbaseStockToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0009965'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '-37.78%'

# Row 49
$ws.Range("B49").Value = 'BOLO'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.002764'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-9.52%'

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002092'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.57%'

# Row 51
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '-0.57%'
